$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.218.45"
$ws.Range("E2").Value = "  +3.77%  "

$ws.Range("D3").Value = "3.080.73"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.36%  "

$ws.Range("D5").Value = "214.99"
$ws.Range("E5").Value = "  +0.50%  "

$ws.Range("D6").Value = "617.65"
$ws.Range("E6").Value = "  -2.58%  "

$ws.Range("D7").Value = "0.375"
$ws.Range("E7").Value = "  -3.04%  "

$ws.Range("D8").Value = "0.875"
$ws.Range("E8").Value = "  +11.58%  "

$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("D10").Value = "3.073.87"
$ws.Range("E10").Value = "  -0.87%  "

$ws.Range("D11").Value = "0.676"
$ws.Range("E11").Value = "  +20.77%  "

$ws.Range("E12").Value = "  +5.73%  "

$ws.Range("D13").Value = "0.0000247"
$ws.Range("E13").Value = "  -0.11%  "

$ws.Range("D14").Value = "90.838.46"
$ws.Range("E14").Value = "  +3.38%  "

$ws.Range("D15").Value = "5.37"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("D16").Value = "32.85"
$ws.Range("E16").Value = "  +2.80%  "

$ws.Range("D17").Value = "3.634.32"
$ws.Range("E17").Value = "  -1.12%  "

$ws.Range("D18").Value = "3.165.85"
$ws.Range("E18").Value = "  +1.46%  "

$ws.Range("D19").Value = "3.46"
$ws.Range("E19").Value = "  +2.79%  "

$ws.Range("D20").Value = "0.0000222"
$ws.Range("E20").Value = "  +2.40%  "

$ws.Range("D21").Value = "13.72"
$ws.Range("E21").Value = "  +4.36%  "

$ws.Range("D22").Value = "432.04"
$ws.Range("E22").Value = "  +2.94%  "

$ws.Range("D23").Value = "8.46"
$ws.Range("E23").Value = "  +1.25%  "

$ws.Range("D24").Value = "5.08"
$ws.Range("E24").Value = "  +4.40%  "

$ws.Range("D25").Value = "5.51"
$ws.Range("E25").Value = "  +1.90%  "

$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "11.81"
$ws.Range("E26").Value = "  +3.29%  "

$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").Value = "83.34"
$ws.Range("E27").Value = "  +1.90%  "

$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("E30").Value = "  +6.51%  "

$ws.Range("D31").Value = "0.166"
$ws.Range("E31").Value = "  +6.67%  "

$ws.Range("D32").Value = "8.66"
$ws.Range("E32").Value = "  +6.80%  "

$ws.Range("D33").Value = "3.83"
$ws.Range("E33").Value = "  -4.17%  "

$ws.Range("D34").Value = "514.15"
$ws.Range("E34").Value = "  +3.02%  "

$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("D36").Value = "1.84"
$ws.Range("E36").Value = "  +0.75%  "

$ws.Range("D37").Value = "22.98"
$ws.Range("E37").Value = "  +4.04%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("E39").Value = "  -7.08%  "

$ws.Range("D40").Value = "22.31"
$ws.Range("E40").Value = "  +0.59%  "

$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "0.366"
$ws.Range("E43").Value = "  +1.08%  "

$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "1.87"
$ws.Range("E44").Value = "  +1.91%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.138"
$ws.Range("E45").Value = "  +3.64%  "

$ws.Range("D46").Value = "0.0719"
$ws.Range("E46").Value = "  +10.27%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "143.35"
$ws.Range("E47").Value = "  -1.67%  "

$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").Value = "0.000271"
$ws.Range("E48").Value = "  +15.49%  "

$ws.Range("D49").Value = "43.73"
$ws.Range("E49").Value = "  +0.36%  "

$ws.Range("D50").Value = "4.24"
$ws.Range("E50").Value = "  +8.33%  "

$ws.Range("D51").Value = "165.03"
$ws.Range("E51").Value = "  +2.43%  "
